{"js": "// The original document has a single paragraph \"123\" followed by the\n// \"_GoBack\" bookmark (an empty, collapsed bookmark Word maintains at the\n// location of the most recent edit). The target state splits that single\n// paragraph into three paragraphs (\"123\", \"456\", \"hgefjhgjhfe\"), with the\n// \"_GoBack\" bookmark now trailing the final paragraph's run - exactly as\n// if the author had clicked at the end of \"123\" and typed\n// Enter, \"456\", Enter, \"hgefjhgjhfe\".\n\n// Word keeps \"_GoBack\" unique, so first remove the existing one; we'll\n// re-insert it at the new final location once the text is in place.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst body = context.document.body;\nconst firstPara = body.paragraphs.getFirst();\n\n// Insert the paragraph breaks + new text right at the end of the \"123\"\n// paragraph's text (i.e. where the bookmark used to sit).\nconst insertionPoint = firstPara.getRange(Word.RangeLocation.end);\ninsertionPoint.insertText(\"\\r456\\rhgefjhgjhfe\", Word.InsertLocation.before);\nawait context.sync();\n\n// Re-create \"_GoBack\" at the very end of the document, after the newly\n// typed \"hgefjhgjhfe\" text, matching Word's normal behavior of tracking\n// the last edit point.\nconst endOfDoc = body.getRange(Word.RangeLocation.end);\nendOfDoc.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Word keeps only one bookmark per name; remove the existing \"_GoBack\" and\n# re-add it at the new final edit location once all the new text is in place.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# --- Split the \"123\" paragraph into three paragraphs: \"123\", \"456\",\n# \"hgefjhgjhfe\" - as if someone clicked at the end of \"123\" and typed\n# Enter, \"456\", Enter, \"hgefjhgjhfe\".\n\n$p1 = $d.Paragraphs.Item(1).Range\n$p1.Collapse(0)              # wdCollapseEnd: end of \"123\"'s text\n$p1.InsertParagraphAfter()   # paragraph 1 (\"123\") / paragraph 2 (new, empty)\n\n$p2 = $d.Paragraphs.Item(2).Range\n$p2.Collapse(1)              # wdCollapseStart: start of the new (empty) paragraph 2\n$p2.InsertAfter(\"456\")\n\n$p2b = $d.Paragraphs.Item(2).Range\n$p2b.Collapse(0)             # end of \"456\"'s text\n$p2b.InsertParagraphAfter()  # paragraph 2 (\"456\") / paragraph 3 (new, empty)\n\n$p3 = $d.Paragraphs.Item(3).Range\n$p3.Collapse(1)              # start of the new (empty) paragraph 3\n$startPos = $p3.Start\n$text3 = \"hgefjhgjhfe\"\n$p3.InsertAfter($text3)\n$anchorPos = $startPos + $text3.Length   # position right after \"hgefjhgjhfe\"\n\n# --- Re-anchor \"_GoBack\" immediately after \"hgefjhgjhfe\".\n#\n# Quirk: adding a bookmark to a collapsed range that sits exactly at the end\n# of a paragraph's stored text (right before its paragraph mark) makes this\n# engine wrap the bookmark around the whole paragraph run instead of\n# collapsing it in place. Work around it by temporarily appending a\n# placeholder character after the anchor point (so the bookmark's position is\n# genuinely \"interior\" to the paragraph, not at its edge), anchoring the\n# bookmark there, then deleting the placeholder.\n$placeholderRange = $d.Range($anchorPos, $anchorPos)\n$placeholderRange.InsertAfter(\"X\")\n\n$bmRange = $d.Range($anchorPos, $anchorPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$placeholder = $d.Range($anchorPos, $anchorPos + 1)\n$placeholder.Delete()\n"}
